$wb = $excel.ActiveWorkbook
$cgs = $wb.Worksheets.Item("CGs")

# Add the new sheet right after "CGs"
$new = $wb.Worksheets.Add($null, $cgs)
$new.Name = "CGs with constellation"

# Header row
$new.Range("A1").Value = 'Names'
$new.Range("B1").Value = 'CG Name'
$new.Range("C1").Value = 'RAJ2000'
$new.Range("D1").Value = 'DEJ2000'
$new.Range("E1").Value = "Constellation"
$new.Range("F1").Value = "Constellation short"

# Row 2
$new.Range("A2").Value = 'CG1 DN / BDN256.11-14.12 / HMSTG256.2-14.1'
$new.Range("B2").Formula = '=IFERROR(MID(A2,1,FIND(" /", A2)),A2)'
$new.Range("C2").Value = '07 19 01'
$new.Range("D2").Value = '-44 34 35'
$new.Range("E2").Value = 'Puppis'
$new.Range("F2").Value = 'Pup'

# Row 3
$new.Range("A3").Value = 'CG2 DN / HMSTG255.3-14.4'
$new.Range("B3:B42").Formula = '=IFERROR(MID(A3,1,FIND(" /", A3)),A3)'
$new.Range("C3").Value = '07 16 00'
$new.Range("D3").Value = '-43 57 41'
$new.Range("E3").Value = 'Puppis'
$new.Range("F3").Value = 'Pup'

# Row 4
$new.Range("A4").Value = 'CG3 DN / FeSt2-35 / HMSTG260.7-12.4'
$new.Range("C4").Value = '07 39 14'
$new.Range("D4").Value = '-47 51 57'
$new.Range("E4").Value = 'Puppis'
$new.Range("F4").Value = 'Pup'

# Row 5
$new.Range("A5").Value = 'CG4 DN / HMSTG259.4-12.7'
$new.Range("C5").Value = '07 34 09'
$new.Range("D5").Value = '-46 54 19'
$new.Range("E5").Value = 'Puppis'
$new.Range("F5").Value = 'Pup'

# Row 6
$new.Range("A6").Value = 'CG5 DN / HMSTG257.2-10.3'
$new.Range("C6").Value = '07 40 51'
$new.Range("D6").Value = '-43 49 09'
$new.Range("E6").Value = 'Puppis'
$new.Range("F6").Value = 'Pup'

# Row 7
$new.Range("A7").Value = 'CG6 DN / FeSt2-27 / HMSTG259.0-13.2'
$new.Range("C7").Value = '07 30 33'
$new.Range("D7").Value = '-46 43 40'
$new.Range("E7").Value = 'Puppis'
$new.Range("F7").Value = 'Pup'

# Row 8
$new.Range("A8").Value = 'CG7 DN / FeSt2-45 / HMSTG266.0-4.3'
$new.Range("C8").Value = '09 14 07'
$new.Range("D8").Value = '-42 30 04'
$new.Range("E8").Value = 'Vela'
$new.Range("F8").Value = 'Vel'

# Row 9
$new.Range("A9").Value = 'CG7S DN'
$new.Range("C9").Value = '02 34 39'
$new.Range("D9").Value = '+61 23 18'
$new.Range("E9").Value = 'Cassiopeia'
$new.Range("F9").Value = 'Cas'

# Row 10
$new.Range("A10").Value = 'CG8 DN / HMSTG255.1-8.8'
$new.Range("C10").Value = '07 42 39'
$new.Range("D10").Value = '-41 15 58'
$new.Range("E10").Value = 'Puppis'
$new.Range("F10").Value = 'Pup'

# Row 11
$new.Range("A11").Value = 'CG9 DN / HMSTG255.1-9.2'
$new.Range("C11").Value = '07 40 47'
$new.Range("D11").Value = '-41 26 57'
$new.Range("E11").Value = 'Puppis'
$new.Range("F11").Value = 'Pup'

# Row 12
$new.Range("A12").Value = 'CG10 DN / HMSTG255.8-9.2'
$new.Range("C12").Value = '07 42 30'
$new.Range("D12").Value = '-42 05 10'
$new.Range("E12").Value = 'Puppis'
$new.Range("F12").Value = 'Pup'

# Row 13
$new.Range("A13").Value = 'CG11 DN / HMSTG4.9-24.6'
$new.Range("C13").Value = '19 40 38'
$new.Range("D13").Value = '-34 47 36'
$new.Range("E13").Value = 'Sagittarius'
$new.Range("F13").Value = 'Sgr'

# Row 14
$new.Range("A14").Value = 'CG12 DN / BDN316.46+21.13 / HMSTG316.5+21.2 /  KM316.5+21.0 rel RN'
$new.Range("C14").Value = '13 57 36'
$new.Range("D14").Value = '-39 59 36'
$new.Range("E14").Value = 'Centaurus'
$new.Range("F14").Value = 'Cen'

# Row 15
$new.Range("A15").Value = 'CG13 DN / HMSTG259.5-16.4 / KM259.5-16.5'
$new.Range("C15").Value = '07 14 12'
$new.Range("D15").Value = '-48 29 10'
$new.Range("E15").Value = 'Puppis'
$new.Range("F15").Value = 'Pup'

# Row 16
$new.Range("A16").Value = 'CG14 DN / FeSt2-37 / HMSTG262.5-13.4'
$new.Range("C16").Value = '07 38 33'
$new.Range("D16").Value = '-49 52 55'
$new.Range("E16").Value = 'Puppis'
$new.Range("F16").Value = 'Pup'

# Row 17
$new.Range("A17").Value = 'CG15 DN / FeSt2-39 / HMSTG262.9-14.7'
$new.Range("C17").Value = '07 32 24'
$new.Range("D17").Value = '-50 46 30'
$new.Range("E17").Value = 'Puppis'
$new.Range("F17").Value = 'Pup'

# Row 18
$new.Range("A18").Value = 'CG16 DN / FeSt2-38 / HMSTG262.9-15.5'
$new.Range("C18").Value = '07 27 47'
$new.Range("D18").Value = '-51 05 11'
$new.Range("E18").Value = 'Carina'
$new.Range("F18").Value = 'Car'

# Row 19
$new.Range("A19").Value = 'CG17 DN / HMSTG270.6-4.7'
$new.Range("C19").Value = '08 52 35'
$new.Range("D19").Value = '-51 51 53'
$new.Range("E19").Value = 'Vela'
$new.Range("F19").Value = 'Vel'

# Row 20
$new.Range("A20").Value = 'CG18 DN / HMSTG269.7-3.9 / VMF24 / V24'
$new.Range("C20").Value = '08 52 38'
$new.Range("D20").Value = '-50 40 05'
$new.Range("E20").Value = 'Vela'
$new.Range("F20").Value = 'Vel'

# Row 21
$new.Range("A21").Value = 'CG19 DN / HMSTG302.1+7.4'
$new.Range("C21").Value = '12 45 40'
$new.Range("D21").Value = '-55 25 23'
$new.Range("E21").Value = 'Centaurus'
$new.Range("F21").Value = 'Cen'

# Row 22
$new.Range("A22").Value = 'CG20 DN / SDN151 / HMSTG302.0-7.0 / BHR82'
$new.Range("C22").Value = '12 40 54'
$new.Range("D22").Value = '-69 52 15'
$new.Range("E22").Value = 'Musca'
$new.Range("F22").Value = 'Mus'

# Row 23
$new.Range("A23").Value = 'CG21 DN / SDN149 / HMSTG301.7-7.2 / BHR80'
$new.Range("C23").Value = '12 37 08'
$new.Range("D23").Value = '-69 59 54'
$new.Range("E23").Value = 'Musca'
$new.Range("F23").Value = 'Mus'

# Row 24
$new.Range("A24").Value = 'CG22 DN / HMSTG253.6+2.9'
$new.Range("C24").Value = '08 28 46'
$new.Range("D24").Value = '-33 45 50'
$new.Range("E24").Value = 'Pyxis'
$new.Range("F24").Value = 'Pyx'

# Row 25
$new.Range("A25").Value = 'CG24 DN / HMSTG260.0-3.8 / VMF21 / V21'
$new.Range("C25").Value = '08 19 15'
$new.Range("D25").Value = '-42 54 47'
$new.Range("E25").Value = 'Puppis'
$new.Range("F25").Value = 'Pup'

# Row 26
$new.Range("A26").Value = 'CG25 DN / HMSTG260.6-12.7'
$new.Range("C26").Value = '07 37 22'
$new.Range("D26").Value = '-47 56 32'
$new.Range("E26").Value = 'Puppis'
$new.Range("F26").Value = 'Pup'

# Row 27
$new.Range("A27").Value = 'CG26 DN / FeSt2-11 / HMSTG252.2+0.7'
$new.Range("C27").Value = '08 15 50'
$new.Range("D27").Value = '-33 50 16'
$new.Range("E27").Value = 'Puppis'
$new.Range("F27").Value = 'Pup'

# Row 28
$new.Range("A28").Value = 'CG27 DN / HMSTG251.7+0.2'
$new.Range("C28").Value = '08 12 25'
$new.Range("D28").Value = '-33 45 39'
$new.Range("E28").Value = 'Puppis'
$new.Range("F28").Value = 'Pup'

# Row 29
$new.Range("A29").Value = 'CG28 DN / FeSt2-9 / HMSTG251.8+0.0 / VMF15 / V15'
$new.Range("C29").Value = '08 12 21'
$new.Range("D29").Value = '-33 56 15'
$new.Range("E29").Value = 'Puppis'
$new.Range("F29").Value = 'Pup'

# Row 30
$new.Range("A30").Value = 'CG29 DN / HMSTG251.9+0.0'
$new.Range("C30").Value = '08 12 22'
$new.Range("D30").Value = '-34 01 03'
$new.Range("E30").Value = 'Puppis'
$new.Range("F30").Value = 'Pup'

# Row 31
$new.Range("A31").Value = 'CG30 DN / HMSTG253.3-1.6 / BHR12 / LM88 / VMF14 / V14'
$new.Range("C31").Value = '08 09 33'
$new.Range("D31").Value = '-36 04 57'
$new.Range("E31").Value = 'Puppis'
$new.Range("F31").Value = 'Pup'

# Row 32
$new.Range("A32").Value = 'CG31 DN / FeSt2-15 / HMSTG253.1-1.7 / VMF13 / V13'
$new.Range("C32").Value = '08 08 56'
$new.Range("D32").Value = '-35 59 21'
$new.Range("E32").Value = 'Puppis'
$new.Range("F32").Value = 'Pup'

# Row 33
$new.Range("A33").Value = 'CG31A DN / BHR8 / LM87 in HMSTG253.1-1.7'
$new.Range("C33").Value = '08 09 02'
$new.Range("D33").Value = '-36 01 10'
$new.Range("E33").Value = 'Puppis'
$new.Range("F33").Value = 'Pup'

# Row 34
$new.Range("A34").Value = 'CG31B DN / BHR9 / LM85 in HMSTG253.1-1.7'
$new.Range("C34").Value = '08 08 48'
$new.Range("D34").Value = '-36 03 12'
$new.Range("E34").Value = 'Puppis'
$new.Range("F34").Value = 'Pup'

# Row 35
$new.Range("A35").Value = 'CG31C DN / BHR10'
$new.Range("C35").Value = '08 08 33'
$new.Range("D35").Value = '-35 59 34'
$new.Range("E35").Value = 'Puppis'
$new.Range("F35").Value = 'Pup'

# Row 36
$new.Range("A36").Value = 'CG31D DN / BHR11'
$new.Range("C36").Value = '08 08 16'
$new.Range("D36").Value = '-36 01 47'
$new.Range("E36").Value = 'Puppis'
$new.Range("F36").Value = 'Pup'

# Row 37
$new.Range("A37").Value = 'CG32 DN / FeSt2-12 / HMSTG252.5+0.1 / VMF18 / V18'
$new.Range("C37").Value = '08 14 21'
$new.Range("D37").Value = '-34 30 53'
$new.Range("E37").Value = 'Puppis'
$new.Range("F37").Value = 'Pup'

# Row 38
$new.Range("A38").Value = 'CG33 DN / HMSTG252.3+0.5 / VMF19 / V19'
$new.Range("C38").Value = '08 15 28'
$new.Range("D38").Value = '-34 04 45'
$new.Range("E38").Value = 'Puppis'
$new.Range("F38").Value = 'Pup'

# Row 39
$new.Range("A39").Value = 'CG34 DN / HMSTG253.8-10.9'
$new.Range("C39").Value = '07 29 32'
$new.Range("D39").Value = '-41 10 30'
$new.Range("E39").Value = 'Puppis'
$new.Range("F39").Value = 'Pup'

# Row 40
$new.Range("A40").Value = 'CG36 DN / HMSTG256.9+2.6'
$new.Range("C40").Value = '08 37 18'
$new.Range("D40").Value = '-36 37 55'
$new.Range("E40").Value = 'Pyxis'
$new.Range("F40").Value = 'Pyx'

# Row 41
$new.Range("A41").Value = 'CG37 DN'
$new.Range("C41").Value = '08 12 28'
$new.Range("D41").Value = '-33 05 35'
$new.Range("E41").Value = 'Puppis'
$new.Range("F41").Value = 'Pup'

# Row 42
$new.Range("A42").Value = 'CG38 DN'
$new.Range("C42").Value = '08 09 39'
$new.Range("D42").Value = '-36 10 35'
$new.Range("E42").Value = 'Puppis'
$new.Range("F42").Value = 'Pup'
